$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 4.9144577670714
$ws.Range("D2").Value = 4.192297529291719
$ws.Range("E2").Value = 9.844811837702748
$ws.Range("F2").Value = 63.19838470746566
$ws.Range("G2").Value = 3.786863320881616
$ws.Range("J2").Value = 10.18219733259556
$ws.Range("K2").Value = 22.26100550071206
$ws.Range("M2").Value = 21.04794975097906
$ws.Range("C3").Value = 4.764170278639956
$ws.Range("D3").Value = 4.097106461735187
$ws.Range("E3").Value = 9.876234483678727
$ws.Range("F3").Value = 62.41309039692587
$ws.Range("G3").Value = 3.792010228982551
$ws.Range("J3").Value = 10.18904741310702
$ws.Range("K3").Value = 22.11860087170103
$ws.Range("M3").Value = 21.05434082282557
$ws.Range("C4").Value = 4.671263840977444
$ws.Range("D4").Value = 4.037857382557633
$ws.Range("E4").Value = 9.896681638100556
$ws.Range("F4").Value = 61.93553163132538
$ws.Range("G4").Value = 3.795327324002663
$ws.Range("J4").Value = 10.19410777217012
$ws.Range("K4").Value = 22.0395462194871
$ws.Range("M4").Value = 21.06474306679536
$ws.Range("C5").Value = 4.633311232819747
$ws.Range("D5").Value = 4.013533531416929
$ws.Range("E5").Value = 9.905305167312338
$ws.Range("F5").Value = 61.74223048899783
$ws.Range("G5").Value = 3.796718698615686
$ws.Range("J5").Value = 10.19638443781098
$ws.Range("K5").Value = 22.00946802254621
$ws.Range("M5").Value = 21.07060864839944
$ws.Range("C6").Value = 4.627005637789988
$ws.Range("D6").Value = 4.009484402879812
$ws.Range("E6").Value = 9.90675471427503
$ws.Range("F6").Value = 61.71021621171528
$ws.Range("G6").Value = 3.796952134355596
$ws.Range("J6").Value = 10.19677542400773
$ws.Range("K6").Value = 22.00460341923398
$ws.Range("M6").Value = 21.07168077531736
$ws.Range("C7").Value = 4.670752283715339
$ws.Range("D7").Value = 4.03753003904763
$ws.Range("E7").Value = 9.896796757732879
$ws.Range("F7").Value = 61.93291921670936
$ws.Range("G7").Value = 3.795345927855018
$ws.Range("J7").Value = 10.19413760787774
$ws.Range("K7").Value = 22.03913188513224
$ws.Range("M7").Value = 21.06481558997431
$ws.Range("C8").Value = 4.86280755402375
$ws.Range("D8").Value = 4.159655008604842
$ws.Range("E8").Value = 9.855407605118979
$ws.Range("F8").Value = 62.92674681937657
$ws.Range("G8").Value = 3.788605527360772
$ws.Range("J8").Value = 10.18438172693164
$ws.Range("K8").Value = 22.21018354678636
$ws.Range("M8").Value = 21.04880772877137
$ws.Range("C9").Value = 5.231743559106902
$ws.Range("D9").Value = 4.391855557343565
$ws.Range("E9").Value = 9.783347440774973
$ws.Range("F9").Value = 64.90565391618163
$ws.Range("G9").Value = 3.776623899503989
$ws.Range("J9").Value = 10.17204408522666
$ws.Range("K9").Value = 22.61063352669225
$ws.Range("M9").Value = 21.06888025457229
$ws.Range("C10").Value = 5.494819549326595
$ws.Range("D10").Value = 4.556842178636916
$ws.Range("E10").Value = 9.735887339398598
$ws.Range("F10").Value = 66.36891080180135
$ws.Range("G10").Value = 3.768562826777345
$ws.Range("J10").Value = 10.16714173746286
$ws.Range("K10").Value = 22.94221726886203
$ws.Range("M10").Value = 21.11502369531758
$ws.Range("C11").Value = 5.6121608135658
$ws.Range("D11").Value = 4.630455192632503
$ws.Range("E11").Value = 9.715472652796532
$ws.Range("F11").Value = 67.03470807033651
$ws.Range("G11").Value = 3.765054168326468
$ws.Range("J11").Value = 10.16581955791879
$ws.Range("K11").Value = 23.10062660048485
$ws.Range("M11").Value = 21.14281498740142
$ws.Range("C12").Value = 5.656215143988625
$ws.Range("D12").Value = 4.658106153950816
$ws.Range("E12").Value = 9.707910010980155
$ws.Range("F12").Value = 67.28669134464197
$ws.Range("G12").Value = 3.763748104195266
$ws.Range("J12").Value = 10.1654497469535
$ws.Range("K12").Value = 23.16165025157262
$ws.Range("M12").Value = 21.15431326202499
$ws.Range("C13").Value = 5.646744832211302
$ws.Range("D13").Value = 4.65216128457031
$ws.Range("E13").Value = 9.709531307595517
$ws.Range("F13").Value = 67.2324310014077
$ws.Range("G13").Value = 3.764028386767784
$ws.Range("J13").Value = 10.16552356707467
$ws.Range("K13").Value = 23.14846241155679
$ws.Range("M13").Value = 21.15179363882195
$ws.Range("C14").Value = 5.615793050984279
$ws.Range("D14").Value = 4.632734653255042
$ws.Range("E14").Value = 9.714847109468179
$ws.Range("F14").Value = 67.05544260403715
$ws.Range("G14").Value = 3.764946265988415
$ws.Range("J14").Value = 10.1657865087837
$ws.Range("K14").Value = 23.10562656105711
$ws.Range("M14").Value = 21.14374145515133
$ws.Range("C15").Value = 5.596783419808369
$ws.Range("D15").Value = 4.620805518845087
$ws.Range("E15").Value = 9.718125034122899
$ws.Range("F15").Value = 66.94700914504652
$ws.Range("G15").Value = 3.765511429762724
$ws.Range("J15").Value = 10.16596462030515
$ws.Range("K15").Value = 23.07952192069068
$ws.Range("M15").Value = 21.1389360193326
$ws.Range("C16").Value = 5.48710014206591
$ws.Range("D16").Value = 4.552000882187774
$ws.Range("E16").Value = 9.737245046199506
$ws.Range("F16").Value = 66.32538859641191
$ws.Range("G16").Value = 3.768795296837184
$ws.Range("J16").Value = 10.16724644165944
$ws.Range("K16").Value = 22.93201286158294
$ws.Range("M16").Value = 21.11334403091237
$ws.Range("C17").Value = 5.419182113285316
$ws.Range("D17").Value = 4.509409681634002
$ws.Range("E17").Value = 9.749274805362573
$ws.Range("F17").Value = 65.94396814310311
$ws.Range("G17").Value = 3.77085027302094
$ws.Range("J17").Value = 10.16826556411054
$ws.Range("K17").Value = 22.84342412193562
$ws.Range("M17").Value = 21.09938373941258
$ws.Range("C18").Value = 5.379900217365511
$ws.Range("D18").Value = 4.484777761768621
$ws.Range("E18").Value = 9.756304676660326
$ws.Range("F18").Value = 65.72460928349662
$ws.Range("G18").Value = 3.772047157955523
$ws.Range("J18").Value = 10.16893717661927
$ws.Range("K18").Value = 22.79318598746571
$ws.Range("M18").Value = 21.09199452227774
$ws.Range("C19").Value = 5.366564168257733
$ws.Range("D19").Value = 4.476415257665369
$ws.Range("E19").Value = 9.758703909721477
$ws.Range("F19").Value = 65.65034725270762
$ws.Range("G19").Value = 3.77245497022797
$ws.Range("J19").Value = 10.16917923781302
$ws.Range("K19").Value = 22.77630067113868
$ws.Range("M19").Value = 21.0896027352786
$ws.Range("C20").Value = 5.426434925534225
$ws.Range("D20").Value = 4.513957648318687
$ws.Range("E20").Value = 9.747982770301959
$ws.Range("F20").Value = 65.98456968807884
$ws.Range("G20").Value = 3.770629974778616
$ws.Range("J20").Value = 10.16814823175157
$ws.Range("K20").Value = 22.85278082857365
$ws.Range("M20").Value = 21.10080358238869
$ws.Range("C21").Value = 5.624894993676854
$ws.Range("D21").Value = 4.638446957751967
$ws.Range("E21").Value = 9.713281179467696
$ws.Range("F21").Value = 67.10743345559848
$ws.Range("G21").Value = 3.764676051028517
$ws.Range("J21").Value = 10.1657057222447
$ws.Range("K21").Value = 23.11818075019472
$ws.Range("M21").Value = 21.14608016813253
$ws.Range("C22").Value = 5.752366074703046
$ws.Range("D22").Value = 4.718490982629085
$ws.Range("E22").Value = 9.691580191831639
$ws.Range("F22").Value = 67.84040244589049
$ws.Range("G22").Value = 3.760916396130461
$ws.Range("J22").Value = 10.16487234564903
$ws.Range("K22").Value = 23.29765697019299
$ws.Range("M22").Value = 21.18134763631038
$ws.Range("C23").Value = 5.684550395206806
$ws.Range("D23").Value = 4.675895966341789
$ws.Range("E23").Value = 9.70307322222336
$ws.Range("F23").Value = 67.44933708042093
$ws.Range("G23").Value = 3.762911016123896
$ws.Range("J23").Value = 10.16524722529628
$ws.Range("K23").Value = 23.20133355213201
$ws.Range("M23").Value = 21.16200678790656
$ws.Range("C24").Value = 5.423156657385951
$ws.Range("D24").Value = 4.511901964965925
$ws.Range("E24").Value = 9.748566544805621
$ws.Range("F24").Value = 65.96621395139634
$ws.Range("G24").Value = 3.77072952347128
$ws.Range("J24").Value = 10.16820101078155
$ws.Range("K24").Value = 22.84854850000971
$ws.Range("M24").Value = 21.10015968781772
$ws.Range("C25").Value = 5.133117146976083
$ws.Range("D25").Value = 4.329944074636518
$ws.Range("E25").Value = 9.801873988706159
$ws.Range("F25").Value = 64.36801376260887
$ws.Range("G25").Value = 3.779734118251278
$ws.Range("J25").Value = 10.17465226873986
$ws.Range("K25").Value = 22.49556035715184
$ws.Range("M25").Value = 21.05793385317333
